# Update gh-pages to output generated at 456a3b4
# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1) updates
$wsExhibit.Range("F3").Value = 13183
$wsExhibit.Range("F6").Value = 104
$wsExhibit.Range("F9").Value = 37
$wsExhibit.Range("F11").Value = 13133
$wsExhibit.Range("F14").Value = 8808
$wsExhibit.Range("F15").Value = 7878
$wsExhibit.Range("F21").Value = 4
$wsExhibit.Range("F27").Value = 73
$wsExhibit.Range("F28").Value = 346

# Sheet "全部类型" (sheet4) updates
$wsAll.Range("F4").Value = 13183
$wsAll.Range("F7").Value = 104
$wsAll.Range("F10").Value = 37
$wsAll.Range("F12").Value = 13133
$wsAll.Range("F15").Value = 8808
$wsAll.Range("F16").Value = 7878
$wsAll.Range("F22").Value = 4
$wsAll.Range("F30").Value = 73
$wsAll.Range("F31").Value = 346
